# Refresh the crypto price ("D") and 1h volume/% change ("E") columns with
# the latest scraped values (commit: "Updated cryptos list on Sun Mar 26
# 11:18:43 UTC 2023 with GitHub Actions").
#
# The sheet stores these as plain text (e.g. "27.761.44", "1.001",
# "  +0.96%  ") even though some look numeric, so each write is forced to
# stay text via a leading quote-prefix and the cell style is reset back to
# "Normal" afterwards to avoid picking up an explicit text number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Addr = "D2"; Val = '27.761.44' },
    @{ Addr = "E2"; Val = '  +0.96%  ' },
    @{ Addr = "D3"; Val = '1.775.30' },
    @{ Addr = "E3"; Val = '  +1.50%  ' },
    @{ Addr = "D4"; Val = '1.001' },
    @{ Addr = "E4"; Val = '  -0.09%  ' },
    @{ Addr = "D5"; Val = '327.35' },
    @{ Addr = "E5"; Val = '  +0.99%  ' },
    @{ Addr = "D6"; Val = '1.000' },
    @{ Addr = "E6"; Val = '  +0.01%  ' },
    @{ Addr = "D7"; Val = '0.4581' },
    @{ Addr = "E7"; Val = '  +2.65%  ' },
    @{ Addr = "D8"; Val = '0.3585' },
    @{ Addr = "E8"; Val = '  -0.38%  ' },
    @{ Addr = "D9"; Val = '0.07497' },
    @{ Addr = "E9"; Val = '  +0.06%  ' },
    @{ Addr = "D10"; Val = '41.87' },
    @{ Addr = "E10"; Val = '  -0.23%  ' },
    @{ Addr = "D11"; Val = '1.107' },
    @{ Addr = "E11"; Val = '  +1.45%  ' },
    @{ Addr = "D12"; Val = '1.001' },
    @{ Addr = "E12"; Val = '  -0.10%  ' },
    @{ Addr = "D13"; Val = '20.85' },
    @{ Addr = "E13"; Val = '  +1.23%  ' },
    @{ Addr = "D14"; Val = '6.051' },
    @{ Addr = "E14"; Val = '  +0.51%  ' },
    @{ Addr = "D15"; Val = '7.230' },
    @{ Addr = "E15"; Val = '  +1.52%  ' },
    @{ Addr = "D16"; Val = '1.771.45' },
    @{ Addr = "E16"; Val = '  +1.09%  ' },
    @{ Addr = "D17"; Val = '93.81' },
    @{ Addr = "E17"; Val = '  +0.83%  ' },
    @{ Addr = "D18"; Val = '0.00001061' },
    @{ Addr = "E18"; Val = '  +0.01%  ' },
    @{ Addr = "D19"; Val = '0.06436' },
    @{ Addr = "E19"; Val = '  +0.86%  ' },
    @{ Addr = "D20"; Val = '1.000' },
    @{ Addr = "E20"; Val = '  +0.03%  ' },
    @{ Addr = "D21"; Val = '17.11' },
    @{ Addr = "E21"; Val = '  +2.11%  ' },
    @{ Addr = "D22"; Val = '5.814' },
    @{ Addr = "E22"; Val = '  -0.66%  ' },
    @{ Addr = "D23"; Val = '27.793.65' },
    @{ Addr = "D24"; Val = '11.33' },
    @{ Addr = "E24"; Val = '  +1.44%  ' },
    @{ Addr = "E25"; Val = '  +0.18%  ' },
    @{ Addr = "D26"; Val = '164.48' },
    @{ Addr = "E26"; Val = '  +1.79%  ' },
    @{ Addr = "D27"; Val = '20.31' },
    @{ Addr = "E27"; Val = '  -0.91%  ' },
    @{ Addr = "D28"; Val = '1.977.48' },
    @{ Addr = "E28"; Val = '  +1.41%  ' },
    @{ Addr = "D29"; Val = '2.182' },
    @{ Addr = "E29"; Val = '  +4.50%  ' },
    @{ Addr = "D30"; Val = '125.78' },
    @{ Addr = "E30"; Val = '  +0.11%  ' },
    @{ Addr = "D31"; Val = '1.105' },
    @{ Addr = "E31"; Val = '  +2.32%  ' },
    @{ Addr = "D32"; Val = '0.09219' },
    @{ Addr = "E32"; Val = '  +2.24%  ' },
    @{ Addr = "D33"; Val = '3.668' },
    @{ Addr = "E33"; Val = '  +0.42%  ' },
    @{ Addr = "D34"; Val = '5.553' },
    @{ Addr = "E34"; Val = '  +0.31%  ' },
    @{ Addr = "D35"; Val = '11.89' },
    @{ Addr = "E35"; Val = '  -0.33%  ' },
    @{ Addr = "D36"; Val = '0.02297' },
    @{ Addr = "D37"; Val = '0.06183' },
    @{ Addr = "E37"; Val = '  +2.83%  ' },
    @{ Addr = "D38"; Val = '0.2092' },
    @{ Addr = "E38"; Val = '  +0.42%  ' },
    @{ Addr = "D39"; Val = '0.6334' },
    @{ Addr = "E39"; Val = '  -0.08%  ' },
    @{ Addr = "D40"; Val = '4.968' },
    @{ Addr = "E40"; Val = '  +0.51%  ' },
    @{ Addr = "D41"; Val = '1.187' },
    @{ Addr = "E41"; Val = '  -1.53%  ' },
    @{ Addr = "D42"; Val = '1.387' },
    @{ Addr = "E42"; Val = '  +0.65%  ' },
    @{ Addr = "D43"; Val = '7.827' },
    @{ Addr = "E43"; Val = '  +1.13%  ' },
    @{ Addr = "D44"; Val = '13.26' },
    @{ Addr = "E44"; Val = '  +0.23%  ' },
    @{ Addr = "D45"; Val = '3.746' },
    @{ Addr = "E45"; Val = '  +0.79%  ' },
    @{ Addr = "D46"; Val = '0.5920' },
    @{ Addr = "E46"; Val = '  +0.65%  ' },
    @{ Addr = "D47"; Val = '122.88' },
    @{ Addr = "E47"; Val = '  +0.66%  ' },
    @{ Addr = "D48"; Val = '1.956' },
    @{ Addr = "E48"; Val = '  +0.16%  ' },
    @{ Addr = "D49"; Val = '0.06926' },
    @{ Addr = "E49"; Val = '  +0.96%  ' },
    @{ Addr = "D50"; Val = '1.139' },
    @{ Addr = "E50"; Val = '  -0.58%  ' },
    @{ Addr = "D51"; Val = '72.51' },
    @{ Addr = "E51"; Val = '  +0.74%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Addr)
    $cell.Value = "'" + $u.Val
    $cell.Style = "Normal"
}
